# TAKEN feedback format and Questions update
#
# - Rename the TASK.0 feedback column header (G1) to the new quiz-based label.
# - Insert a new student row (Maaya Leonard / c1243957) as row 4, pushing the
#   existing roster rows down.
# - Refresh the per-row "Last.Access" (E), "Group Code" (L), mirrored
#   First/Last name (N/O) and "Group Set" (P) values for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the feedback column header text -------------------------------
$ws.Cells.Item(1, 7).Value = "TASK.0..data...quiz.links..Total.Pts..1.Score...114592"

# --- Insert the new roster row for Maaya Leonard (c1243957) ---------------
$ws.Rows.Item(4).Insert()

# --- Row 2: c1206235 / Vandebroek / Martina --------------------------------
$ws.Cells.Item(2, 1).Value = "c1206235"
$ws.Cells.Item(2, 2).Value = "Vandebroek"
$ws.Cells.Item(2, 3).Value = "Martina "
$ws.Cells.Item(2, 4).Value = 31206235
$ws.Cells.Item(2, 5).Value = "'"
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(2, 6).Value = "Yes"
$ws.Cells.Item(2, 12).Value = "ATSTAT"
$ws.Cells.Item(2, 13).Value = 31206235
$ws.Cells.Item(2, 14).Value = "Martina "
$ws.Cells.Item(2, 15).Value = "Vandebroek"
$ws.Cells.Item(2, 16).Value = "3_gc_groups"
$ws.Cells.Item(2, 17).Value = "ePAxyD"

# --- Row 3: c1242115 / Vandebroek / Martina --------------------------------
$ws.Cells.Item(3, 1).Value = "c1242115"
$ws.Cells.Item(3, 2).Value = "Vandebroek"
$ws.Cells.Item(3, 3).Value = "Martina"
$ws.Cells.Item(3, 4).Value = 31242115
$ws.Cells.Item(3, 5).Value = "'"
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(3, 6).Value = "Yes"
$ws.Cells.Item(3, 12).Value = "TSTAT"
$ws.Cells.Item(3, 13).Value = 31242115
$ws.Cells.Item(3, 14).Value = "Martina"
$ws.Cells.Item(3, 15).Value = "Vandebroek"
$ws.Cells.Item(3, 16).Value = "3_gc_groups"
$ws.Cells.Item(3, 17).Value = "iVQYPJ"

# --- Row 4 (new): c1243957 / Maaya / Leonard -------------------------------
$ws.Cells.Item(4, 1).Value = "c1243957"
$ws.Cells.Item(4, 2).Value = "Maaya"
$ws.Cells.Item(4, 3).Value = "Leonard"
$ws.Cells.Item(4, 5).Value = "2023-08-23 13:54:33"
$ws.Cells.Item(4, 6).Value = "Yes"
$ws.Cells.Item(4, 12).Value = "TSTAT"
$ws.Cells.Item(4, 14).Value = "Leonard"
$ws.Cells.Item(4, 15).Value = "Maaya"
$ws.Cells.Item(4, 16).Value = "3_gc_groups"
$ws.Cells.Item(4, 17).Value = "yXzkle"

# --- Row 5: q0762379 / Assele / Samson Yaekob ------------------------------
$ws.Cells.Item(5, 1).Value = "q0762379"
$ws.Cells.Item(5, 2).Value = "Assele"
$ws.Cells.Item(5, 3).Value = "Samson Yaekob"
$ws.Cells.Item(5, 4).Value = 445993
$ws.Cells.Item(5, 5).Value = "'"
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(5, 6).Value = "Yes"
$ws.Cells.Item(5, 12).Value = "ATSTAT"
$ws.Cells.Item(5, 13).Value = 445993
$ws.Cells.Item(5, 14).Value = "Samson Yaekob"
$ws.Cells.Item(5, 15).Value = "Assele"
$ws.Cells.Item(5, 16).Value = "3_gc_groups"
$ws.Cells.Item(5, 17).Value = "hPQFvR"

# --- Row 6: q1371623 / Gutierrez Vargas / Alvaro ---------------------------
$ws.Cells.Item(6, 1).Value = "q1371623"
$ws.Cells.Item(6, 2).Value = "Gutierrez Vargas"
$ws.Cells.Item(6, 3).Value = "Alvaro"
$ws.Cells.Item(6, 4).Value = 50133260
$ws.Cells.Item(6, 5).Value = "'"
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(6, 6).Value = "Yes"
$ws.Cells.Item(6, 12).Value = "MMENG"
$ws.Cells.Item(6, 13).Value = 50133260
$ws.Cells.Item(6, 14).Value = "Alvaro"
$ws.Cells.Item(6, 15).Value = "Gutierrez Vargas"
$ws.Cells.Item(6, 16).Value = "3_gc_groups"
$ws.Cells.Item(6, 17).Value = "lDLHdt"

# --- Row 7: q1411379 / Xi / Yuan -------------------------------------------
$ws.Cells.Item(7, 1).Value = "q1411379"
$ws.Cells.Item(7, 2).Value = "Xi"
$ws.Cells.Item(7, 3).Value = "Yuan"
$ws.Cells.Item(7, 4).Value = 818343
$ws.Cells.Item(7, 5).Value = "2023-08-23 10:46:27"
$ws.Cells.Item(7, 6).Value = "Yes"
$ws.Cells.Item(7, 12).Value = "MMENG"
$ws.Cells.Item(7, 13).Value = 818343
$ws.Cells.Item(7, 14).Value = "Yuan"
$ws.Cells.Item(7, 15).Value = "Xi"
$ws.Cells.Item(7, 16).Value = "3_gc_groups"
$ws.Cells.Item(7, 17).Value = "PrlRAD"
